$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.597.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.89%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.855.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.07%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.16%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'522.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +6.01%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'141.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.41%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.609"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.22%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.713"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.69%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -5.05%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0000321"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -8.15%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'41.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.71%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'10.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.07%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.468.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.26%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'21.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +7.69%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.866.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.26%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'14.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.38%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -2.03%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +1.43%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'68.575.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.01%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'415.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.53%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.21%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'14.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.79%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'12.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.10%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'86.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.28%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'4.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +6.25%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -5.98%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'35.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.40%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'13.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.12%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'675.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.41%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +15.47%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -4.38%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'2.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.50%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'66.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +8.26%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.452"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.55%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.0₃0849"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -6.77%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'39.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.70%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +11.89%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -2.42%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.19%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.24%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.30%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0476"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.48%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'3.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +5.52%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.65%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -1.45%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.000283"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +17.13%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'3.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.55%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'3.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.74%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'8.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.80%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'142.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.53%  "
$ws.Range("E51").Style = "Normal"
